$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.998.48"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.642.86"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D5").Value = "'216.37"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").Value = "'19.66"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").Value = "'0.0796"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "1.868.98"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "1.659.69"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "'63.08"
$ws.Range("D18").Value = "25.974.49"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'193.30"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'9.94"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'6.26"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +7.08%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "'144.52"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("E35").Value = "  +2.94%  "
$ws.Range("D36").Value = "'0.906"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "1.133.15"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'0.541"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").Value = "'99.23"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'0.797"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "1.778.21"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("D46").Value = "'56.69"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "'0.0529"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +0.33%  "
